$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the score for the Fortuna United vs Ofside match (row 8)
$ws.Range("F8").Value = 6
$ws.Range("G8").Value = 2

# Update the active cell selection to F9
$ws.Range("F9").Select()
